$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.086.50"
$ws.Range("E2").Value = "  -1.25%  "

$ws.Range("D3").Value = "2.467.32"
$ws.Range("E3").Value = "  -2.88%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.69"
$ws.Range("E5").Value = "  -1.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.92"
$ws.Range("E6").Value = "  -2.52%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -2.25%  "

$ws.Range("D9").Value = "2.465.27"
$ws.Range("E9").Value = "  -2.95%  "

$ws.Range("E10").Value = "  -2.84%  "

$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.94"
$ws.Range("E12").Value = "  -2.82%  "

$ws.Range("E13").Value = "  -3.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.57"
$ws.Range("E14").Value = "  -3.43%  "

$ws.Range("D16").Value = "67.061.01"
$ws.Range("E16").Value = "  -1.18%  "

$ws.Range("E17").Value = "  -4.75%  "

$ws.Range("D18").Value = "2.436.55"
$ws.Range("E18").Value = "  -3.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.11"
$ws.Range("E19").Value = "  -5.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  -4.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.96"
$ws.Range("E21").Value = "  -4.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.04"
$ws.Range("E22").Value = "  -2.85%  "

$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.94"
$ws.Range("E24").Value = "  -3.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.26"
$ws.Range("E25").Value = "  -7.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.80"
$ws.Range("E26").Value = "  -6.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.23"
$ws.Range("E27").Value = "  -7.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -59.22%  "

$ws.Range("D29").Value = "2.590.64"
$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0904"
$ws.Range("E30").Value = "  -6.88%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "515.06"
$ws.Range("E31").Value = "  -4.84%  "

$ws.Range("E32").Value = "  -8.47%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.24"
$ws.Range("E33").Value = "  -6.14%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.77"
$ws.Range("E34").Value = "  -5.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").Value = "  -9.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.73"
$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.35"
$ws.Range("E39").Value = "  -4.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  -6.31%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.82"
$ws.Range("E42").Value = "  -6.60%  "

$ws.Range("E43").Value = "  -6.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.327"
$ws.Range("E44").Value = "  -7.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").Value = "  -7.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.78"
$ws.Range("E46").Value = "  -1.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.65"
$ws.Range("E47").Value = "  -4.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.46"
$ws.Range("E48").Value = "  -6.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.515"
$ws.Range("E49").Value = "  -6.92%  "

$ws.Range("E50").Value = "  -11.50%  "

$ws.Range("E51").Value = "  -7.04%  "

